# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh values to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4728.593
$ws.Range("I100").Value = 1780.3636
$ws.Range("J100").Value = 6755.5
$ws.Range("K100").Value = 1780.3636
$ws.Range("L100").Value = 6755.5
$ws.Range("M100").Value = -1239.3636
$ws.Range("N100").Value = -7837.5

$ws.Range("H113").Value = 4989.9375
$ws.Range("I113").Value = 2779.8572
$ws.Range("J113").Value = 6708.8887
$ws.Range("K113").Value = 2779.8572
$ws.Range("L113").Value = 6708.8887
$ws.Range("M113").Value = 474.1428000000001
$ws.Range("N113").Value = -13216.8887

$ws.Range("H116").Value = 98725.55
$ws.Range("I116").Value = 108100.1
$ws.Range("J116").Value = 4980
$ws.Range("K116").Value = 108100.1
$ws.Range("L116").Value = 4980
$ws.Range("M116").Value = -104658.1
$ws.Range("N116").Value = -11864

$ws.Range("H135").Value = 400.125
$ws.Range("I135").Value = 226.8
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 2041.2
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = 493.8
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1190.8334
$ws.Range("I2").Value = 1266.2858
$ws.Range("J2").Value = 662.6667
$ws.Range("K2").Value = 1266.2858
$ws.Range("L2").Value = 662.6667
$ws.Range("M2").Value = -1153.2858
$ws.Range("N2").Value = -888.6667

$ws.Range("H32").Value = 1211375.9
$ws.Range("I32").Value = 1519156.8
$ws.Range("J32").Value = 2236.7144
$ws.Range("K32").Value = 1519156.8
$ws.Range("L32").Value = 2236.7144
$ws.Range("M32").Value = -1518869.8
$ws.Range("N32").Value = -2810.7144

$ws.Range("H45").Value = 980.5599999999999
$ws.Range("I45").Value = 915.65
$ws.Range("J45").Value = 1240.2
$ws.Range("K45").Value = 915.65
$ws.Range("L45").Value = 1240.2
$ws.Range("M45").Value = -538.65
$ws.Range("N45").Value = -1994.2

$ws.Range("H88").Value = 3622.3777
$ws.Range("I88").Value = 3817.9487
$ws.Range("J88").Value = 2351.1667
$ws.Range("K88").Value = 3817.9487
$ws.Range("L88").Value = 2351.1667
$ws.Range("M88").Value = -3411.9487
$ws.Range("N88").Value = -3163.1667

$ws.Range("H91").Value = 3622.3777
$ws.Range("I91").Value = 3817.9487
$ws.Range("J91").Value = 2351.1667
$ws.Range("K91").Value = 3817.9487
$ws.Range("L91").Value = 2351.1667
$ws.Range("M91").Value = -2413.9487
$ws.Range("N91").Value = -5159.1667

$ws.Range("H109").Value = 200000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 200000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 200000
$ws.Range("N109").Value = -202774

$ws.Range("H110").Value = 1937.9474
$ws.Range("I110").Value = 1951.3125
$ws.Range("J110").Value = 1866.6666
$ws.Range("K110").Value = 1951.3125
$ws.Range("L110").Value = 1866.6666
$ws.Range("M110").Value = 93.6875
$ws.Range("N110").Value = -5956.6666

$ws.Range("H116").Value = 1190.8334
$ws.Range("I116").Value = 1266.2858
$ws.Range("J116").Value = 662.6667
$ws.Range("K116").Value = 1266.2858
$ws.Range("L116").Value = 662.6667
$ws.Range("M116").Value = 1027.7142
$ws.Range("N116").Value = -5250.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1190.8334
$ws.Range("I3").Value = 1266.2858
$ws.Range("J3").Value = 662.6667
$ws.Range("K3").Value = 1266.2858
$ws.Range("L3").Value = 662.6667
$ws.Range("M3").Value = -1152.2858
$ws.Range("N3").Value = -890.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H31").Value = 2658.8367
$ws.Range("I31").Value = 1371.0646
$ws.Range("J31").Value = 4876.6665
$ws.Range("K31").Value = 1371.0646
$ws.Range("L31").Value = 4876.6665
$ws.Range("M31").Value = -1076.0646
$ws.Range("N31").Value = -5466.6665

$ws.Range("H34").Value = 2658.8367
$ws.Range("I34").Value = 1371.0646
$ws.Range("J34").Value = 4876.6665
$ws.Range("K34").Value = 1371.0646
$ws.Range("L34").Value = 4876.6665
$ws.Range("M34").Value = -1169.0646
$ws.Range("N34").Value = -5280.6665

$ws.Range("H58").Value = 5395.5713
$ws.Range("I58").Value = 9133
$ws.Range("J58").Value = 2592.5
$ws.Range("K58").Value = 9133
$ws.Range("L58").Value = 2592.5
$ws.Range("M58").Value = -8930
$ws.Range("N58").Value = -2998.5

$ws.Range("H99").Value = 55031.156
$ws.Range("I99").Value = 145384.28
$ws.Range("J99").Value = 2325.1667
$ws.Range("K99").Value = 145384.28
$ws.Range("L99").Value = 2325.1667
$ws.Range("M99").Value = -143886.28
$ws.Range("N99").Value = -5321.1667

$ws.Range("H126").Value = 55031.156
$ws.Range("I126").Value = 145384.28
$ws.Range("J126").Value = 2325.1667
$ws.Range("K126").Value = 436152.84
$ws.Range("L126").Value = 6975.500100000001
$ws.Range("M126").Value = -433682.84
$ws.Range("N126").Value = -11915.5001

$ws.Range("H134").Value = 1744.1666
$ws.Range("I134").Value = 996.1429000000001
$ws.Range("J134").Value = 2791.4
$ws.Range("K134").Value = 2988.4287
$ws.Range("L134").Value = 8374.200000000001
$ws.Range("M134").Value = -453.4287000000004
$ws.Range("N134").Value = -13444.2

$ws.Range("H136").Value = 5395.5713
$ws.Range("I136").Value = 9133
$ws.Range("J136").Value = 2592.5
$ws.Range("K136").Value = 27399
$ws.Range("L136").Value = 7777.5
$ws.Range("M136").Value = -24849
$ws.Range("N136").Value = -12877.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 7451.6
$ws.Range("I94").Value = 3024
$ws.Range("J94").Value = 7943.5557
$ws.Range("K94").Value = 9072
$ws.Range("L94").Value = 23830.6671
$ws.Range("M94").Value = -8396
$ws.Range("N94").Value = -25182.6671

$ws.Range("H101").Value = 4857.143
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 4857.143
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 14571.429
$ws.Range("N101").Value = -19439.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2898.4
$ws.Range("I7").Value = 2830.6667
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2830.6667
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2718.6667
$ws.Range("N7").Value = -3224

$ws.Range("H8").Value = 2898.4
$ws.Range("I8").Value = 2830.6667
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 2830.6667
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -2691.6667
$ws.Range("N8").Value = -3278

$ws.Range("H12").Value = 25001.334
$ws.Range("I12").Value = 2500
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -2360
$ws.Range("N12").Value = -70284

$ws.Range("H126").Value = 2078.9524
$ws.Range("I126").Value = 1688.4706
$ws.Range("J126").Value = 3738.5
$ws.Range("K126").Value = 5065.4118
$ws.Range("L126").Value = 11215.5
$ws.Range("M126").Value = -2595.4118
$ws.Range("N126").Value = -16155.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H40").Value = 2162.1333
$ws.Range("I40").Value = 2144.3333
$ws.Range("J40").Value = 2233.3333
$ws.Range("K40").Value = 2144.3333
$ws.Range("L40").Value = 2233.3333
$ws.Range("M40").Value = -2008.3333
$ws.Range("N40").Value = -2505.3333

$ws.Range("H61").Value = 1359
$ws.Range("I61").Value = 1440.8
$ws.Range("J61").Value = 950
$ws.Range("K61").Value = 1440.8
$ws.Range("L61").Value = 950
$ws.Range("M61").Value = -1238.8
$ws.Range("N61").Value = -1354

$ws.Range("H113").Value = 1359
$ws.Range("I113").Value = 1440.8
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 1440.8
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 729.2
$ws.Range("N113").Value = -5290

$ws.Range("H132").Value = 11236.444
$ws.Range("I132").Value = 4776.231
$ws.Range("J132").Value = 17235.215
$ws.Range("K132").Value = 14328.693
$ws.Range("L132").Value = 51705.645
$ws.Range("M132").Value = -11798.693
$ws.Range("N132").Value = -56765.645

$ws.Range("H136").Value = 5325.1333
$ws.Range("I136").Value = 3013.7896
$ws.Range("J136").Value = 9317.454
$ws.Range("K136").Value = 9041.3688
$ws.Range("L136").Value = 27952.362
$ws.Range("M136").Value = -6491.3688
$ws.Range("N136").Value = -33052.362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 2450
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2450
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2450
$ws.Range("N15").Value = -3026

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H107").Value = 408.82144
$ws.Range("I107").Value = 322.95
$ws.Range("J107").Value = 623.5
$ws.Range("K107").Value = 968.8499999999999
$ws.Range("L107").Value = 1870.5
$ws.Range("M107").Value = 951.1500000000001
$ws.Range("N107").Value = -5710.5

$ws.Range("H126").Value = 1657.7
$ws.Range("I126").Value = 1125
$ws.Range("J126").Value = 3788.5
$ws.Range("K126").Value = 3375
$ws.Range("L126").Value = 11365.5
$ws.Range("M126").Value = -905
$ws.Range("N126").Value = -16305.5

$ws.Range("H132").Value = 2747.6775
$ws.Range("I132").Value = 2054.0908
$ws.Range("J132").Value = 4443.1113
$ws.Range("K132").Value = 6162.2724
$ws.Range("L132").Value = 13329.3339
$ws.Range("M132").Value = -3632.2724
$ws.Range("N132").Value = -18389.3339

$ws.Range("H136").Value = 18724186
$ws.Range("I136").Value = 23280530
$ws.Range("J136").Value = 913027.75
$ws.Range("K136").Value = 69841590
$ws.Range("L136").Value = 2739083.25
$ws.Range("M136").Value = -69839040
$ws.Range("N136").Value = -2744183.25

